$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 183071.81
$ws.Range("I28").Value = 222421.22
$ws.Range("J28").Value = 5999.5
$ws.Range("K28").Value = 222421.22
$ws.Range("L28").Value = 5999.5
$ws.Range("M28").Value = -221936.22
$ws.Range("N28").Value = -6969.5

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5874.269
$ws.Range("I62").Value = 4991.1
$ws.Range("J62").Value = 8818.166999999999
$ws.Range("K62").Value = 4991.1
$ws.Range("L62").Value = 8818.166999999999
$ws.Range("M62").Value = -4367.1
$ws.Range("N62").Value = -10066.167

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5874.269
$ws.Range("I65").Value = 4991.1
$ws.Range("J65").Value = 8818.166999999999
$ws.Range("K65").Value = 24955.5
$ws.Range("L65").Value = 44090.835
$ws.Range("M65").Value = -21835.5
$ws.Range("N65").Value = -50330.835

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 26316462
$ws.Range("I92").Value = 26316462
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 26316462
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -26315214

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3155.361
$ws.Range("I100").Value = 2209.2856
$ws.Range("K100").Value = 2209.2856
$ws.Range("M100").Value = -1668.2856

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5735.7856
$ws.Range("I113").Value = 4261.8
$ws.Range("K113").Value = 4261.8
$ws.Range("M113").Value = -1007.8

# ALC row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1245
$ws.Range("I115").Value = 614
$ws.Range("J115").Value = 4400
$ws.Range("K115").Value = 1842
$ws.Range("L115").Value = 13200
$ws.Range("M115").Value = -275
$ws.Range("N115").Value = -16334

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1053.421
$ws.Range("I135").Value = 950.8333
$ws.Range("K135").Value = 8557.4997
$ws.Range("M135").Value = -6022.4997

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1378204.6
$ws.Range("I137").Value = 70574.336
$ws.Range("J137").Value = 1919293
$ws.Range("K137").Value = 211723.008
$ws.Range("L137").Value = 5757879
$ws.Range("M137").Value = -209173.008
$ws.Range("N137").Value = -5762979

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4086.78
$ws.Range("I138").Value = 1802.2858
$ws.Range("J138").Value = 4458.6743
$ws.Range("K138").Value = 5406.857400000001
$ws.Range("L138").Value = 13376.0229
$ws.Range("M138").Value = -266.8574000000008
$ws.Range("N138").Value = -23656.0229

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1168.65
$ws.Range("I2").Value = 1007.35297
$ws.Range("J2").Value = 2082.6667
$ws.Range("K2").Value = 1007.35297
$ws.Range("L2").Value = 2082.6667
$ws.Range("M2").Value = -894.35297
$ws.Range("N2").Value = -2308.6667

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16324448
$ws.Range("I32").Value = 16758730
$ws.Range("J32").Value = 11909252
$ws.Range("K32").Value = 16758730
$ws.Range("L32").Value = 11909252
$ws.Range("M32").Value = -16758443
$ws.Range("N32").Value = -11909826

# ARM row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 51372.9
$ws.Range("I37").Value = 25267
$ws.Range("J37").Value = 68776.836
$ws.Range("K37").Value = 25267
$ws.Range("L37").Value = 68776.836
$ws.Range("M37").Value = -24994
$ws.Range("N37").Value = -69322.836

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 64005.285
$ws.Range("I55").Value = 30048
$ws.Range("J55").Value = 69664.836
$ws.Range("K55").Value = 30048
$ws.Range("L55").Value = 69664.836
$ws.Range("M55").Value = -29733
$ws.Range("N55").Value = -70294.836

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4511.1035
$ws.Range("I61").Value = 4420.143
$ws.Range("J61").Value = 4749.875
$ws.Range("K61").Value = 4420.143
$ws.Range("L61").Value = 4749.875
$ws.Range("M61").Value = -4208.143
$ws.Range("N61").Value = -5173.875

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3192.3076
$ws.Range("I74").Value = 3125
$ws.Range("K74").Value = 3125
$ws.Range("M74").Value = -2251

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3192.3076
$ws.Range("I77").Value = 3125
$ws.Range("K77").Value = 15625
$ws.Range("M77").Value = -11257

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1168.65
$ws.Range("I116").Value = 1007.35297
$ws.Range("J116").Value = 2082.6667
$ws.Range("K116").Value = 1007.35297
$ws.Range("L116").Value = 2082.6667
$ws.Range("M116").Value = 1286.64703
$ws.Range("N116").Value = -6670.6667

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4511.1035
$ws.Range("I136").Value = 4420.143
$ws.Range("J136").Value = 4749.875
$ws.Range("K136").Value = 13260.429
$ws.Range("L136").Value = 14249.625
$ws.Range("M136").Value = -10710.429
$ws.Range("N136").Value = -19349.625

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1168.65
$ws.Range("I3").Value = 1007.35297
$ws.Range("J3").Value = 2082.6667
$ws.Range("K3").Value = 1007.35297
$ws.Range("L3").Value = 2082.6667
$ws.Range("M3").Value = -893.35297
$ws.Range("N3").Value = -2310.6667

# BSM row 62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50181
$ws.Range("J62").Value = 50181
$ws.Range("L62").Value = 50181
$ws.Range("N62").Value = -51553

# BSM row 63
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 84556.42999999999
$ws.Range("J63").Value = 84556.42999999999
$ws.Range("L63").Value = 84556.42999999999
$ws.Range("N63").Value = -85928.42999999999

# BSM row 65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 50181
$ws.Range("J65").Value = 50181
$ws.Range("L65").Value = 150543
$ws.Range("N65").Value = -157407

# BSM row 66
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H66").Value = 84556.42999999999
$ws.Range("J66").Value = 84556.42999999999
$ws.Range("L66").Value = 253669.29
$ws.Range("N66").Value = -260533.29

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3410.2593
$ws.Range("I107").Value = 3428.375
$ws.Range("K107").Value = 3428.375
$ws.Range("M107").Value = -1508.375

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3075.7856
$ws.Range("I31").Value = 1627.1666
$ws.Range("K31").Value = 1627.1666
$ws.Range("M31").Value = -1332.1666

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3075.7856
$ws.Range("I34").Value = 1627.1666
$ws.Range("K34").Value = 1627.1666
$ws.Range("M34").Value = -1425.1666

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2005.7273
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 2056.3
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 2056.3
$ws.Range("M94").Value = -1049
$ws.Range("N94").Value = -2958.3

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1975.8214
$ws.Range("I134").Value = 1916.5217
$ws.Range("K134").Value = 5749.5651
$ws.Range("M134").Value = -3214.5651

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 182675340
$ws.Range("I4").Value = 148345710
$ws.Range("J4").Value = 319993900
$ws.Range("K4").Value = 445037130
$ws.Range("L4").Value = 959981700
$ws.Range("M4").Value = -445037018
$ws.Range("N4").Value = -959981924

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1149.6666
$ws.Range("J5").Value = 1179.8
$ws.Range("L5").Value = 3539.4
$ws.Range("N5").Value = -3763.4

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 672466.7
$ws.Range("I11").Value = 721671.4399999999
$ws.Range("J11").Value = 500250
$ws.Range("K11").Value = 2165014.32
$ws.Range("L11").Value = 1500750
$ws.Range("M11").Value = -2164874.32
$ws.Range("N11").Value = -1501030

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3631.4211
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 4499.625
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 13498.875
$ws.Range("M68").Value = -8189
$ws.Range("N68").Value = -15120.875

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3631.4211
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 4499.625
$ws.Range("K71").Value = 27000
$ws.Range("L71").Value = 40496.625
$ws.Range("M71").Value = -22944
$ws.Range("N71").Value = -48608.625

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 769.4
$ws.Range("I107").Value = 769.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2308.2
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -388.1999999999998

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 479176
$ws.Range("J132").Value = 912381.4399999999
$ws.Range("L132").Value = 8211432.959999999
$ws.Range("N132").Value = -8216492.959999999

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1149.6666
$ws.Range("J135").Value = 1179.8
$ws.Range("L135").Value = 10618.2
$ws.Range("N135").Value = -15688.2

# GSM row 7
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 36673336
$ws.Range("I7").Value = 10010001
$ws.Range("J7").Value = 90000000
$ws.Range("K7").Value = 10010001
$ws.Range("L7").Value = 90000000
$ws.Range("M7").Value = -10009889
$ws.Range("N7").Value = -90000224

# GSM row 8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 36673336
$ws.Range("I8").Value = 10010001
$ws.Range("J8").Value = 90000000
$ws.Range("K8").Value = 10010001
$ws.Range("L8").Value = 90000000
$ws.Range("M8").Value = -10009862
$ws.Range("N8").Value = -90000278

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6770.6895
$ws.Range("I46").Value = 2399.8333
$ws.Range("J46").Value = 7910.913
$ws.Range("K46").Value = 2399.8333
$ws.Range("L46").Value = 7910.913
$ws.Range("M46").Value = -2211.8333
$ws.Range("N46").Value = -8286.913

# WVR row 98
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 100000
$ws.Range("J98").Value = 100000
$ws.Range("L98").Value = 100000
$ws.Range("N98").Value = -105990
